$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 32
$ws1.Range("F8").Value = 3770
$ws1.Range("F9").Value = 76
$ws1.Range("F10").Value = 4438

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 32
$ws4.Range("F9").Value = 3770
$ws4.Range("F10").Value = 76
$ws4.Range("F11").Value = 4438
